
$d = $word.ActiveDocument

# --- Hunk 1: insert new centered paragraph "«Основные конструкции языка Python»." ---
$target = $d.Content
$target.Find.Execute("по курсу «Разработка интернет приложений».", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $target.Paragraphs(1)
$para.Range.InsertParagraphAfter()
$newParaIndex = $para.Index + 1
$p1 = $d.Paragraphs($newParaIndex)
$p1.Range.Text = "«"
$r1 = $p1.Range.Duplicate
$r1.Collapse(0)
$r1.InsertParagraphAfter()
$p2 = $d.Paragraphs($newParaIndex + 1)
$p2.Range.Text = "Основные конструкции языка Python"
$r2 = $p2.Range.Duplicate
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$p3 = $d.Paragraphs($newParaIndex + 2)
$p3.Range.Text = "»."

$markStart = $d.Paragraphs($newParaIndex).Range.End - 1
$markRange = $d.Range($markStart, $markStart + 1)
$markRange.Delete()

$markStart2 = $d.Paragraphs($newParaIndex).Range.End - 1
$markRange2 = $d.Range($markStart2, $markStart2 + 1)
$markRange2.Delete()

# --- Hunk 2: remove one of the 5 duplicate empty centered paragraphs before "Москва, 202" ---
$r = $d.Content
$r.Find.Execute("Москва, 202", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$moscowPara = $r.Paragraphs(1)
$prevPara = $moscowPara.Previous()
$prevPara.Range.Delete()

# --- Hunk 3: merge the run-fragmented sentence about coefficients into a single run ---
$oldText = "Если коэффициент А, В, С введён или задан в командной строке некорректно, то необходимо проигнорировать некорректное значение и вводить коэффициент повторно пока коэффициент не будет введён корректно. Корректно заданный коэффициент — это коэффициент, значение которого может быть без ошибок преобразовано в действительное число."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $oldText, 2)

# --- Hunk 4: mark the picture-containing runs (InlineShapes) as NoProofing ---
$shapes = $d.InlineShapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $s = $shapes.Item($i)
    $sr = $s.Range
    if ($sr.NoProofing -eq 0) {
        $sr.NoProofing = -1
    }
}

Write-Output "all edits applied"
